$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1723.8077
$ws.Range("J112").Value = 1826.2084
$ws.Range("L112").Value = 5478.6252
$ws.Range("N112").Value = -7694.6252

$ws.Range("H137").Value = 2173.0613
$ws.Range("I137").Value = 1541.8182
$ws.Range("K137").Value = 4625.4546
$ws.Range("M137").Value = -2075.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9062.951999999999
$ws.Range("I61").Value = 7094.6206
$ws.Range("J61").Value = 13453.846
$ws.Range("K61").Value = 7094.6206
$ws.Range("L61").Value = 13453.846
$ws.Range("M61").Value = -6882.6206
$ws.Range("N61").Value = -13877.846

$ws.Range("H74").Value = 5634.7856
$ws.Range("I74").Value = 2709.647
$ws.Range("K74").Value = 2709.647
$ws.Range("M74").Value = -1835.647

$ws.Range("H77").Value = 5634.7856
$ws.Range("I77").Value = 2709.647
$ws.Range("K77").Value = 13548.235
$ws.Range("M77").Value = -9180.235000000001

$ws.Range("H136").Value = 9062.951999999999
$ws.Range("I136").Value = 7094.6206
$ws.Range("J136").Value = 13453.846
$ws.Range("K136").Value = 21283.8618
$ws.Range("L136").Value = 40361.538
$ws.Range("M136").Value = -18733.8618
$ws.Range("N136").Value = -45461.538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2643.75
$ws.Range("I86").Value = 2608.3333
$ws.Range("J86").Value = 2750
$ws.Range("K86").Value = 2608.3333
$ws.Range("L86").Value = 2750
$ws.Range("M86").Value = -1485.3333
$ws.Range("N86").Value = -4996

$ws.Range("H89").Value = 2643.75
$ws.Range("I89").Value = 2608.3333
$ws.Range("J89").Value = 2750
$ws.Range("K89").Value = 13041.6665
$ws.Range("L89").Value = 13750
$ws.Range("M89").Value = -7425.666499999999
$ws.Range("N89").Value = -24982

$ws.Range("H107").Value = 1169.0714
$ws.Range("I107").Value = 782
$ws.Range("J107").Value = 1767.2727
$ws.Range("K107").Value = 782
$ws.Range("L107").Value = 1767.2727
$ws.Range("M107").Value = 1138
$ws.Range("N107").Value = -5607.2727

$ws.Range("H134").Value = 54149.75
$ws.Range("I134").Value = 4581.75
$ws.Range("J134").Value = 128501.75
$ws.Range("K134").Value = 13745.25
$ws.Range("L134").Value = 385505.25
$ws.Range("M134").Value = -11210.25
$ws.Range("N134").Value = -390575.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5598.5625
$ws.Range("I31").Value = 5906.696
$ws.Range("J31").Value = 4811.1113
$ws.Range("K31").Value = 5906.696
$ws.Range("L31").Value = 4811.1113
$ws.Range("M31").Value = -5611.696
$ws.Range("N31").Value = -5401.1113

$ws.Range("H34").Value = 5598.5625
$ws.Range("I34").Value = 5906.696
$ws.Range("J34").Value = 4811.1113
$ws.Range("K34").Value = 5906.696
$ws.Range("L34").Value = 4811.1113
$ws.Range("M34").Value = -5704.696
$ws.Range("N34").Value = -5215.1113

$ws.Range("H58").Value = 1625641.1
$ws.Range("I58").Value = 2599019.2
$ws.Range("J58").Value = 3344.0952
$ws.Range("K58").Value = 2599019.2
$ws.Range("L58").Value = 3344.0952
$ws.Range("M58").Value = -2598816.2
$ws.Range("N58").Value = -3750.0952

$ws.Range("H132").Value = 4827.1
$ws.Range("I132").Value = 5210.148
$ws.Range("J132").Value = 4031.5386
$ws.Range("K132").Value = 15630.444
$ws.Range("L132").Value = 12094.6158
$ws.Range("M132").Value = -13100.444
$ws.Range("N132").Value = -17154.6158

$ws.Range("H134").Value = 2780.4827
$ws.Range("I134").Value = 1897.931
$ws.Range("J134").Value = 3663.0344
$ws.Range("K134").Value = 5693.793
$ws.Range("L134").Value = 10989.1032
$ws.Range("M134").Value = -3158.793
$ws.Range("N134").Value = -16059.1032

$ws.Range("H136").Value = 1625641.1
$ws.Range("I136").Value = 2599019.2
$ws.Range("J136").Value = 3344.0952
$ws.Range("K136").Value = 7797057.600000001
$ws.Range("L136").Value = 10032.2856
$ws.Range("M136").Value = -7794507.600000001
$ws.Range("N136").Value = -15132.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 4361.619
$ws.Range("I123").Value = 2950
$ws.Range("J123").Value = 4510.2104
$ws.Range("K123").Value = 8850
$ws.Range("L123").Value = 13530.6312
$ws.Range("M123").Value = -6400
$ws.Range("N123").Value = -18430.6312

$ws.Range("H129").Value = 2152.8333
$ws.Range("I129").Value = 2696.6667
$ws.Range("J129").Value = 1826.5333
$ws.Range("K129").Value = 8090.000100000001
$ws.Range("L129").Value = 5479.5999
$ws.Range("M129").Value = -3090.000100000001
$ws.Range("N129").Value = -15479.5999

$ws.Range("H130").Value = 5499.6
$ws.Range("I130").Value = 1750
$ws.Range("J130").Value = 7999.3335
$ws.Range("K130").Value = 5250
$ws.Range("L130").Value = 23998.0005
$ws.Range("M130").Value = -230
$ws.Range("N130").Value = -34038.00049999999

$ws.Range("H131").Value = 30573.152
$ws.Range("I131").Value = 3244
$ws.Range("J131").Value = 35453.355
$ws.Range("K131").Value = 9732
$ws.Range("L131").Value = 106360.065
$ws.Range("M131").Value = -4692
$ws.Range("N131").Value = -116440.065

$ws.Range("H133").Value = 4744.706
$ws.Range("I133").Value = 3970
$ws.Range("J133").Value = 4983.077
$ws.Range("K133").Value = 11910
$ws.Range("L133").Value = 14949.231
$ws.Range("M133").Value = -6850
$ws.Range("N133").Value = -25069.231

$ws.Range("H134").Value = 4588.931
$ws.Range("I134").Value = 4313.6875
$ws.Range("J134").Value = 4927.6924
$ws.Range("K134").Value = 12941.0625
$ws.Range("L134").Value = 14783.0772
$ws.Range("M134").Value = -7871.0625
$ws.Range("N134").Value = -24923.0772

$ws.Range("H136").Value = 3094.6
$ws.Range("I136").Value = 1018.3333
$ws.Range("J136").Value = 3984.4285
$ws.Range("K136").Value = 3054.9999
$ws.Range("L136").Value = 11953.2855
$ws.Range("M136").Value = 2045.0001
$ws.Range("N136").Value = -22153.2855

$ws.Range("H137").Value = 33592.65
$ws.Range("I137").Value = 1886.6666
$ws.Range("J137").Value = 69261.875
$ws.Range("K137").Value = 5659.9998
$ws.Range("L137").Value = 207785.625
$ws.Range("M137").Value = -559.9997999999996
$ws.Range("N137").Value = -217985.625

$ws.Range("H138").Value = 7905.579
$ws.Range("I138").Value = 13592.5
$ws.Range("J138").Value = 3769.6365
$ws.Range("K138").Value = 40777.5
$ws.Range("L138").Value = 11308.9095
$ws.Range("M138").Value = -35637.5
$ws.Range("N138").Value = -21588.9095

$ws.Range("H139").Value = 1500807.8
$ws.Range("I139").Value = 2610640.8
$ws.Range("J139").Value = 2533.2
$ws.Range("K139").Value = 7831922.399999999
$ws.Range("L139").Value = 7599.599999999999
$ws.Range("M139").Value = -7826782.399999999
$ws.Range("N139").Value = -17879.6

$ws.Range("H140").Value = 2380.8838
$ws.Range("I140").Value = 1779.3549
$ws.Range("J140").Value = 3934.8333
$ws.Range("K140").Value = 5338.0647
$ws.Range("L140").Value = 11804.4999
$ws.Range("M140").Value = -158.0646999999999
$ws.Range("N140").Value = -22164.4999

$ws.Range("H141").Value = 3389.1365
$ws.Range("I141").Value = 2662
$ws.Range("J141").Value = 4947.2856
$ws.Range("K141").Value = 7986
$ws.Range("L141").Value = 14841.8568
$ws.Range("M141").Value = -2806
$ws.Range("N141").Value = -25201.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5024.1577
$ws.Range("I132").Value = 5069.107
$ws.Range("J132").Value = 4898.3
$ws.Range("K132").Value = 15207.321
$ws.Range("L132").Value = 14694.9
$ws.Range("M132").Value = -12677.321
$ws.Range("N132").Value = -19754.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3335.0188
$ws.Range("I132").Value = 3528.2058
$ws.Range("K132").Value = 10584.6174
$ws.Range("M132").Value = -8054.617400000001

$ws.Range("H136").Value = 4737.2935
$ws.Range("I136").Value = 3563.5957
$ws.Range("J136").Value = 6707.4287
$ws.Range("K136").Value = 10690.7871
$ws.Range("L136").Value = 20122.2861
$ws.Range("M136").Value = -8140.7871
$ws.Range("N136").Value = -25222.2861

